# Daily attendance processing - 2025-11-12 15:26:07
#
# Normalizes the "Recorded By" (column G) values on the active sheet:
#   - "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   - "backup@backdoor.com, system, System" -> "backup@backdoor.com, System, system"
#
# Only rows matching those exact patterns are touched; every other
# "Recorded By" combination (e.g. "backup@backdoor.com, System",
# "admin@admin.com, System", single-author entries, etc.) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($null -eq $text -or $text -eq "") {
        continue
    }

    if ($text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($text -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
